$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '30.364.26'
$ws.Range('E2').Value = '  +0.06%  '

# Row 3
$ws.Range('D3').Value = '1.937.14'
$ws.Range('E3').Value = '  +0.01%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.006'
$ws.Range('E4').Value = '  +0.58%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7641'
$ws.Range('E5').Value = '  +5.50%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '245.73'
$ws.Range('E6').Value = '  -2.55%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.005'
$ws.Range('E7').Value = '  +0.49%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3198'
$ws.Range('E8').Value = '  -3.47%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '27.56'
$ws.Range('E9').Value = '  -1.74%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07003'
$ws.Range('E10').Value = '  -3.65%  '

# Row 11
$ws.Range('E11').Value = '  -3.51%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08009'
$ws.Range('E12').Value = '  -1.20%  '

# Row 13
$ws.Range('D13').Value = '1.933.20'
$ws.Range('E13').Value = '  -0.22%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.354'
$ws.Range('E14').Value = '  -2.40%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '94.48'
$ws.Range('E15').Value = '  -0.48%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '14.45'
$ws.Range('E16').Value = '  -4.51%  '

# Row 17
$ws.Range('D17').Value = '30.366.45'
$ws.Range('E17').Value = '  +0.06%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '254.59'
$ws.Range('E18').Value = '  +0.52%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000007926'
$ws.Range('E19').Value = '  -3.97%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '5.780'

# Row 21
$ws.Range('D21').Value = '2.185.65'
$ws.Range('E21').Value = '  -0.06%  '

# Row 22
$ws.Range('E22').Value = '  +0.23%  '

# Row 23
$ws.Range('E23').Value = '  +0.45%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.688'
$ws.Range('E24').Value = '  -3.96%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.514'
$ws.Range('E25').Value = '  -2.63%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '165.01'
$ws.Range('E26').Value = '  -0.64%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '19.04'
$ws.Range('E27').Value = '  -1.63%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.1336'
$ws.Range('E28').Value = '  +2.79%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.265'
$ws.Range('E29').Value = '  -3.71%  '

# Row 30
$ws.Range('E30').Value = '  +0.91%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.514'
$ws.Range('E31').Value = '  -1.97%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.397'
$ws.Range('E32').Value = '  -1.08%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.113'
$ws.Range('E33').Value = '  -2.42%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.05152'
$ws.Range('E34').Value = '  -1.96%  '

# Row 35
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.275'
$ws.Range('E35').Value = '  +0.49%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7453'
$ws.Range('E36').Value = '  -0.89%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.789'
$ws.Range('E37').Value = '  +0.70%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01950'
$ws.Range('E38').Value = '  -1.20%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.809'
$ws.Range('E39').Value = '  +0.20%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '78.73'
$ws.Range('E40').Value = '  -0.88%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.410'
$ws.Range('E41').Value = '  -0.68%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.4482'
$ws.Range('E42').Value = '  -1.60%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.964'
$ws.Range('E43').Value = '  -3.55%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.004'
$ws.Range('E44').Value = '  +0.29%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.8341'
$ws.Range('E45').Value = '  -1.29%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '101.11'
$ws.Range('E46').Value = '  -0.87%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.800'
$ws.Range('E47').Value = '  -0.20%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.465'
$ws.Range('E48').Value = '  +0.16%  '

# Row 49
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '37.19'
$ws.Range('E49').Value = '  +1.02%  '

# Row 50
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '972.24'
$ws.Range('E50').Value = '  +9.44%  '

# Row 51
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.06022'
$ws.Range('E51').Value = '  -0.45%  '
